$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (case with 380 kV), columns B-F and I-N
$blockBF = New-Object 'object[,]' 24,5
$blockBF[0,0] = 1.02
$blockBF[0,1] = 1.024220495633489
$blockBF[0,2] = 1.047369449335628
$blockBF[0,3] = 1.036075524516319
$blockBF[0,4] = 1.050354430087737
$blockBF[1,0] = 1.02
$blockBF[1,1] = 1.025136526813097
$blockBF[1,2] = 1.047974574763446
$blockBF[1,3] = 1.036862142263826
$blockBF[1,4] = 1.051162928113707
$blockBF[2,0] = 1.02
$blockBF[2,1] = 1.02573012897742
$blockBF[2,2] = 1.048365929656911
$blockBF[2,3] = 1.037372118745484
$blockBF[2,4] = 1.051686689115026
$blockBF[3,0] = 1.02
$blockBF[3,1] = 1.025979886019798
$blockBF[3,2] = 1.048530404816759
$blockBF[3,3] = 1.037586746263034
$blockBF[3,4] = 1.051907021910766
$blockBF[4,0] = 1.02
$blockBF[4,1] = 1.026021833409502
$blockBF[4,2] = 1.048558017899662
$blockBF[4,3] = 1.037622796804117
$blockBF[4,4] = 1.051944025121751
$blockBF[5,0] = 1.02
$blockBF[5,1] = 1.025733465432685
$blockBF[5,2] = 1.048368127582948
$blockBF[5,3] = 1.037374985693942
$blockBF[5,4] = 1.051689632650694
$blockBF[6,0] = 1.02
$blockBF[6,1] = 1.02452989200729
$blockBF[6,2] = 1.047573994787314
$blockBF[6,3] = 1.03634116125265
$blockBF[6,4] = 1.050627538506251
$blockBF[7,0] = 1.02
$blockBF[7,1] = 1.02241575230068
$blockBF[7,2] = 1.046173181001817
$blockBF[7,3] = 1.034527029706569
$blockBF[7,4] = 1.048760745834855
$blockBF[8,0] = 1.02
$blockBF[8,1] = 1.021010913774744
$blockBF[8,2] = 1.04523845765275
$blockBF[8,3] = 1.033322821344802
$blockBF[8,4] = 1.047519535467408
$blockBF[9,0] = 1.02
$blockBF[9,1] = 1.020403707268813
$blockBF[9,2] = 1.04483353806125
$blockBF[9,3] = 1.032802643999915
$blockBF[9,4] = 1.046982890513243
$blockBF[10,0] = 1.02
$blockBF[10,1] = 1.020178329932715
$blockBF[10,2] = 1.044683108297633
$blockBF[10,3] = 1.032609616957764
$blockBF[10,4] = 1.046783680060887
$blockBF[11,0] = 1.02
$blockBF[11,1] = 1.020226666595401
$blockBF[11,2] = 1.044715377070826
$blockBF[11,3] = 1.032651013288062
$blockBF[11,4] = 1.04682640576258
$blockBF[12,0] = 1.02
$blockBF[12,1] = 1.020385074110262
$blockBF[12,2] = 1.044821103979301
$blockBF[12,3] = 1.032786684438579
$blockBF[12,4] = 1.046966421186184
$blockBF[13,0] = 1.02
$blockBF[13,1] = 1.020482696268579
$blockBF[13,2] = 1.044886242622504
$blockBF[13,3] = 1.03287030113963
$blockBF[13,4] = 1.047052705720862
$blockBF[14,0] = 1.02
$blockBF[14,1] = 1.021051235308279
$blockBF[14,2] = 1.045265327263635
$blockBF[14,3] = 1.033357370347991
$blockBF[14,4] = 1.047555168013175
$blockBF[15,0] = 1.02
$blockBF[15,1] = 1.02140815963745
$blockBF[15,2] = 1.045503070869407
$blockBF[15,3] = 1.033663232654636
$blockBF[15,4] = 1.047870566951803
$blockBF[16,0] = 1.02
$blockBF[16,1] = 1.021616453401589
$blockBF[16,2] = 1.045641725351761
$blockBF[16,3] = 1.033841757724369
$blockBF[16,4] = 1.048054611572008
$blockBF[17,0] = 1.02
$blockBF[17,1] = 1.021687494106787
$blockBF[17,2] = 1.04568899993953
$blockBF[17,3] = 1.033902650589191
$blockBF[17,4] = 1.048117379156531
$blockBF[18,0] = 1.02
$blockBF[18,1] = 1.021369854078048
$blockBF[18,2] = 1.045477565008388
$blockBF[18,3] = 1.03363040402018
$blockBF[18,4] = 1.047836719572918
$blockBF[19,0] = 1.02
$blockBF[19,1] = 1.020338422447389
$blockBF[19,2] = 1.044789970701518
$blockBF[19,3] = 1.032746727407613
$blockBF[19,4] = 1.046925186711156
$blockBF[20,0] = 1.02
$blockBF[20,1] = 1.019690882629772
$blockBF[20,2] = 1.044357511624496
$blockBF[20,3] = 1.032192224912438
$blockBF[20,4] = 1.04635278467868
$blockBF[21,0] = 1.02
$blockBF[21,1] = 1.020034064113473
$blockBF[21,2] = 1.044586779027196
$blockBF[21,3] = 1.0324860723151
$blockBF[21,4] = 1.046656157370258
$blockBF[22,0] = 1.02
$blockBF[22,1] = 1.021387162384589
$blockBF[22,2] = 1.045489090062936
$blockBF[22,3] = 1.033645237494079
$blockBF[22,4] = 1.047852013506488
$blockBF[23,0] = 1.02
$blockBF[23,1] = 1.022961505035907
$blockBF[23,2] = 1.046535483639906
$blockBF[23,3] = 1.034995114994272
$blockBF[23,4] = 1.049242780612423

$blockIN = New-Object 'object[,]' 24,6
$blockIN[0,0] = 1.036126349546883
$blockIN[0,1] = 1.029396328622976
$blockIN[0,2] = 1.050132278945543
$blockIN[0,3] = 1.038870303632534
$blockIN[0,4] = 1.053108939929659
$blockIN[0,5] = 1.013815668085104
$blockIN[1,0] = 1.036215098220011
$blockIN[1,1] = 1.029951346433641
$blockIN[1,2] = 1.050549770988679
$blockIN[1,3] = 1.0394665088316
$blockIN[1,4] = 1.053729877734805
$blockIN[1,5] = 1.01400271746571
$blockIN[2,0] = 1.036271080155638
$blockIN[2,1] = 1.030310739536851
$blockIN[2,2] = 1.050819019231176
$blockIN[2,3] = 1.039852618440906
$blockIN[2,4] = 1.054131620500899
$blockIN[2,5] = 1.014123750893179
$blockIN[3,0] = 1.036294268379665
$blockIN[3,1] = 1.030461889375661
$blockIN[3,2] = 1.050931994752889
$blockIN[3,3] = 1.040015015342837
$blockIN[3,4] = 1.054300500514517
$blockIN[3,5] = 1.014174632887109
$blockIN[4,0] = 1.036298141437746
$blockIN[4,1] = 1.030487271641401
$blockIN[4,2] = 1.050950951098775
$blockIN[4,3] = 1.040042286938609
$blockIN[4,4] = 1.054328855419495
$blockIN[4,5] = 1.014183176157728
$blockIN[5,0] = 1.036271391360648
$blockIN[5,1] = 1.030312758970221
$blockIN[5,2] = 1.050820529668225
$blockIN[5,3] = 1.039854788097705
$blockIN[5,4] = 1.054133877135802
$blockIN[5,5] = 1.014124430783243
$blockIN[6,0] = 1.036156641171446
$blockIN[6,1] = 1.029583845081462
$blockIN[6,2] = 1.050273557056885
$blockIN[6,3] = 1.039071725853402
$blockIN[6,4] = 1.053318797141323
$blockIN[6,5] = 1.013878882033419
$blockIN[7,0] = 1.035943409309401
$blockIN[7,1] = 1.028301443243498
$blockIN[7,2] = 1.049302925817407
$blockIN[7,3] = 1.037694420001558
$blockIN[7,4] = 1.051882250713129
$blockIN[7,5] = 1.013446213512272
$blockIN[8,0] = 1.035793887617614
$blockIN[8,1] = 1.027447946322044
$blockIN[8,2] = 1.048651369817931
$blockIN[8,3] = 1.036778012512206
$blockIN[8,4] = 1.050924474759034
$blockIN[8,5] = 1.013157807985004
$blockIN[9,0] = 1.035727405590162
$blockIN[9,1] = 1.027078727930038
$blockIN[9,2] = 1.04836820019257
$blockIN[9,3] = 1.036381641298862
$blockIN[9,4] = 1.050509750555352
$blockIN[9,5] = 1.013032940381207
$blockIN[10,0] = 1.035702450846722
$blockIN[10,1] = 1.026941637710541
$blockIN[10,2] = 1.048262863598462
$blockIN[10,3] = 1.036234478765264
$blockIN[10,4] = 1.050355705149884
$blockIN[10,5] = 1.012986561532699
$blockIN[11,0] = 1.035707815493172
$blockIN[11,1] = 1.026971041569259
$blockIN[11,2] = 1.048285465632694
$blockIN[11,3] = 1.036266042556763
$blockIN[11,4] = 1.050388748310245
$blockIN[11,5] = 1.012996509830976
$blockIN[12,0] = 1.035725348130016
$blockIN[12,1] = 1.027067394902985
$blockIN[12,2] = 1.048359496184647
$blockIN[12,3] = 1.036369475415305
$blockIN[12,4] = 1.05049701707065
$blockIN[12,5] = 1.013029106636707
$blockIN[13,0] = 1.035736116089493
$blockIN[13,1] = 1.027126768552297
$blockIN[13,2] = 1.048405088390944
$blockIN[13,3] = 1.036433212783617
$blockIN[13,4] = 1.050563725295886
$blockIN[13,5] = 1.013049190954548
$blockIN[14,0] = 1.035798263256543
$blockIN[14,1] = 1.027472457631398
$blockIN[14,2] = 1.048670141053577
$blockIN[14,3] = 1.036804327737682
$blockIN[14,4] = 1.050951998770322
$blockIN[14,5] = 1.013166095367916
$blockIN[15,0] = 1.03583678173107
$blockIN[15,1] = 1.027689394140847
$blockIN[15,2] = 1.048836123926489
$blockIN[15,3] = 1.03703723692981
$blockIN[15,4] = 1.051195553468996
$blockIN[15,5] = 1.01323943045645
$blockIN[16,0] = 1.035859081147314
$blockIN[16,1] = 1.027815963365993
$blockIN[16,2] = 1.04893283831721
$blockIN[16,3] = 1.037173131211898
$blockIN[16,4] = 1.051337614599542
$blockIN[16,5] = 1.01328220688299
$blockIN[17,0] = 1.035866656193563
$blockIN[17,1] = 1.027859125904455
$blockIN[17,2] = 1.048965798308056
$blockIN[17,3] = 1.037219474781327
$blockIN[17,4] = 1.051386053685485
$blockIN[17,5] = 1.013296792747466
$blockIN[18,0] = 1.035832666409686
$blockIN[18,1] = 1.027666115399106
$blockIN[18,2] = 1.048818325929452
$blockIN[18,3] = 1.037012243593368
$blockIN[18,4] = 1.051169422363225
$blockIN[18,5] = 1.01323156215671
$blockIN[19,0] = 1.035720192386372
$blockIN[19,1] = 1.027039019749868
$blockIN[19,2] = 1.048337700291527
$blockIN[19,3] = 1.036339015145209
$blockIN[19,4] = 1.050465134560106
$blockIN[19,5] = 1.013019507616147
$blockIN[20,0] = 1.035647969535412
$blockIN[20,1] = 1.026645052415599
$blockIN[20,2] = 1.048034617163797
$blockIN[20,3] = 1.035916120245906
$blockIN[20,4] = 1.050022331163784
$blockIN[20,5] = 1.012886195425268
$blockIN[21,0] = 1.035686398723901
$blockIN[21,1] = 1.026853871944136
$blockIN[21,2] = 1.048195371506946
$blockIN[21,3] = 1.036140267351073
$blockIN[21,4] = 1.050257068057711
$blockIN[21,5] = 1.012956865175277
$blockIN[22,0] = 1.03583452646471
$blockIN[22,1] = 1.027676633956076
$blockIN[22,2] = 1.048826368389273
$blockIN[22,3] = 1.037023536875769
$blockIN[22,4] = 1.051181229885957
$blockIN[22,5] = 1.01323511749868
$blockIN[23,0] = 1.035999836802701
$blockIN[23,1] = 1.028632726424805
$blockIN[23,2] = 1.049554652204297
$blockIN[23,3] = 1.038050175952985
$blockIN[23,4] = 1.052253653819785
$blockIN[23,5] = 1.013558063586324

$ws.Range("B2:F25").Value = $blockBF
$ws.Range("I2:N25").Value = $blockIN

Write-Host "Updated vm_pu values for 380 kV case"
